$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '65.531.49'
$ws.Range('E2').Value = '  -3.36%  '

# Row 3
$ws.Range('D3').Value = '3.492.84'
$ws.Range('E3').Value = '  -0.60%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.02%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '553.50'
$ws.Range('E5').Value = '  -0.39%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.61'
$ws.Range('E6').Value = '  -6.05%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.642'
$ws.Range('E7').Value = '  +4.92%  '

# Row 8
$ws.Range('E8').Value = '  +0.04%  '

# Row 9
$ws.Range('E9').Value = '  -1.09%  '

# Row 10
$ws.Range('E10').Value = '  +2.65%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.02'
$ws.Range('E11').Value = '  -5.42%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.18'
$ws.Range('E13').Value = '  -3.26%  '

# Row 14
$ws.Range('D14').Value = '4.050.75'
$ws.Range('E14').Value = '  -0.71%  '

# Row 15
$ws.Range('D15').Value = '3.489.63'
$ws.Range('E15').Value = '  -0.86%  '

# Row 16
$ws.Range('E16').Value = '  +0.17%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.43'
$ws.Range('E17').Value = '  +0.30%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.21'
$ws.Range('E18').Value = '  +2.46%  '

# Row 19
$ws.Range('D19').Value = '65.519.32'
$ws.Range('E19').Value = '  -4.13%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.995'
$ws.Range('E20').Value = '  -1.38%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '414.80'
$ws.Range('E21').Value = '  +1.20%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '85.75'
$ws.Range('E23').Value = '  +1.10%  '

# Row 24
$ws.Range('E24').Value = '  -2.63%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.78'
$ws.Range('E25').Value = '  +7.27%  '

# Row 26
$ws.Range('E26').Value = '  -7.78%  '

# Row 27
$ws.Range('E27').Value = '  -1.92%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.04'
$ws.Range('E28').Value = '  +4.65%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '30.36'
$ws.Range('E29').Value = '  -0.79%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '621.57'
$ws.Range('E30').Value = '  -8.88%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.47'
$ws.Range('E31').Value = '  -6.14%  '

# Row 32
$ws.Range('E32').Value = '  -0.88%  '

# Row 33
$ws.Range('E33').Value = '  -1.31%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '59.52'
$ws.Range('E34').Value = '  -1.70%  '

# Row 35
$ws.Range('E35').Value = '  +10.45%  '

# Row 36
$ws.Range('E36').Value = '  +0.17%  '

# Row 37
$ws.Range('D37').Value = '0.0₃0791'
$ws.Range('E37').Value = '  -5.88%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '37.12'
$ws.Range('E38').Value = '  -5.11%  '

# Row 39
$ws.Range('D39').Value = '3.370.20'
$ws.Range('E39').Value = '  +10.90%  '

# Row 40
$ws.Range('E40').Value = '  -5.89%  '

# Row 41
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.28'
$ws.Range('E41').Value = '  -4.06%  '

# Row 42
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  -0.11%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.86'
$ws.Range('E43').Value = '  -5.72%  '

# Row 44
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.27'
$ws.Range('E44').Value = '  +0.46%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.52'
$ws.Range('E45').Value = '  -8.49%  '

# Row 46
$ws.Range('E46').Value = '  -2.09%  '

# Row 47
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.77'
$ws.Range('E47').Value = '  +0.33%  '

# Row 48
$ws.Range('E48').Value = '  +1.77%  '

# Row 49
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.42'
$ws.Range('E49').Value = '  -9.96%  '

# Row 50
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '137.31'
$ws.Range('E50').Value = '  -0.87%  '

# Row 51
$ws.Range('E51').Value = '  +9.99%  '
